# Update "想去人数" (number of people interested) figures for several
# events across the "展览" and "全部类型" sheets, plus the single entry
# on the "演出" sheet, matching the regenerated site data.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Row = 2;  Value = 234 },
    @{ Row = 3;  Value = 262 },
    @{ Row = 4;  Value = 276 },
    @{ Row = 6;  Value = 265 },
    @{ Row = 7;  Value = 6473 },
    @{ Row = 11; Value = 75 },
    @{ Row = 14; Value = 3 },
    @{ Row = 15; Value = 205 },
    @{ Row = 16; Value = 516 }
)

# 展览 sheet
$wsExhibit = $wb.Worksheets.Item("展览")
foreach ($u in $updates) {
    $wsExhibit.Cells.Item($u.Row, 6).Value = $u.Value
}

# 演出 sheet
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Cells.Item(2, 6).Value = 6

# 全部类型 sheet (combines both of the above, plus the 演出 row at 18)
$wsAll = $wb.Worksheets.Item("全部类型")
foreach ($u in $updates) {
    $wsAll.Cells.Item($u.Row, 6).Value = $u.Value
}
$wsAll.Cells.Item(18, 6).Value = 6
